$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics for Reln-Itgb1 (rows 2-10, cols E-T)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.041827
$ws.Range("H2").Value = 0.125481
$ws.Range("I2").Value = 0.006279874897961605
$ws.Range("J2").Value = 0.006279874897961606
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 2.553187268104334
$ws.Range("R2").Value = 22.978685412939
$ws.Range("S2").Value = 0.001283363687219847
$ws.Range("T2").Value = 0.001283363687219847
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.041827
$ws.Range("H3").Value = 0.125481
$ws.Range("I3").Value = 0.006279874897961605
$ws.Range("J3").Value = 0.006279874897961606
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 4.446815169382
$ws.Range("R3").Value = 40.02133652443801
$ws.Range("S3").Value = 0.002235198797775778
$ws.Range("T3").Value = 0.002235198797775778
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.041827
$ws.Range("H4").Value = 0.125481
$ws.Range("I4").Value = 0.006279874897961605
$ws.Range("J4").Value = 0.006279874897961606
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 5.493491647185333
$ws.Range("R4").Value = 49.441424824668
$ws.Range("S4").Value = 0.00276131241296598
$ws.Range("T4").Value = 0.002761312412965981
$ws.Range("I5").Value = 0.1693441751896972
$ws.Range("J5").Value = 0.1693441751896972
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 68.84968236585533
$ws.Range("R5").Value = 619.6471412926981
$ws.Range("S5").Value = 0.03460740358875573
$ws.Range("T5").Value = 0.03460740358875573
$ws.Range("I6").Value = 0.1693441751896972
$ws.Range("J6").Value = 0.1693441751896972
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.06027475116060124
$ws.Range("T6").Value = 0.06027475116060126
$ws.Range("I7").Value = 0.1693441751896972
$ws.Range("J7").Value = 0.1693441751896972
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 148.1384306247973
$ws.Range("R7").Value = 1333.245875623176
$ws.Range("S7").Value = 0.07446202044034021
$ws.Range("T7").Value = 0.07446202044034023
$ws.Range("G8").Value = 5.490742
$ws.Range("H8").Value = 16.472226
$ws.Range("I8").Value = 0.8243759499123412
$ws.Range("J8").Value = 0.8243759499123412
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 335.1637116418994
$ws.Range("R8").Value = 3016.473404777094
$ws.Range("S8").Value = 0.1684705787814779
$ws.Range("T8").Value = 0.1684705787814779
$ws.Range("G9").Value = 5.490742
$ws.Range("H9").Value = 16.472226
$ws.Range("I9").Value = 0.8243759499123412
$ws.Range("J9").Value = 0.8243759499123412
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 583.745303673772
$ws.Range("R9").Value = 5253.707733063948
$ws.Range("S9").Value = 0.2934205158700593
$ws.Range("T9").Value = 0.2934205158700593
$ws.Range("G10").Value = 5.490742
$ws.Range("H10").Value = 16.472226
$ws.Range("I10").Value = 0.8243759499123412
$ws.Range("J10").Value = 0.8243759499123412
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 721.1453203397252
$ws.Range("R10").Value = 6490.307883057528
$ws.Range("S10").Value = 0.362484855260804
$ws.Range("T10").Value = 0.3624848552608041
